# Added two new Mac-Addresses
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Two new rows of machine/Mac-Address data, continuing the existing pattern
$newRows = @(
    @(10001, 110030, 10030),
    @(10001, 110031, 10031)
)

$r = 31
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r = $r + 1
}

# Scroll the view down to the newly added rows and select the last-used cell
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("F30").Select()
